$d = $word.ActiveDocument

# --- "Hirse" paragraph --------------------------------------------------
# Original: "Hirse 400g 0,124 mg"  ->  "Hirse 400 g 0,124 mg"
# (a space is inserted between "400" and "g"; formatting/highlight stays green)
$d.Content.Find.Execute("Hirse 400g 0,124 mg", $false, $false, $false, $false, $false, $true, 1, $false, "Hirse 400 g 0,124 mg", 2) | Out-Null

# --- "Quinoa" paragraph --------------------------------------------------
# Original: "Quinoa ?" (highlighted red) followed by a plain trailing space.
# New data was gathered, so the placeholder text is replaced with the real
# figure and the whole line is re-highlighted green (matching the rest of
# the table).
$d.Content.Find.Execute("Quinoa ? ", $false, $false, $false, $false, $false, $true, 1, $false, "Quinoa 390 g 0,075 mg", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Quinoa 390 g 0,075 mg") | Out-Null
$rng.HighlightColorIndex = 4
